# Applies the "Add files via upload" revision:
#   - Slide 15 "TextBox 7": bump both bullet runs to 20pt and grow the
#     (spAutoFit) text box to its new rendered height.
#   - Slide 4 "Content Placeholder 2": fix "Analyzing" -> "Analyze" in two runs.
#   - Slide 6 table (Content Placeholder 5): header cell "Methodology" ->
#     "Project Methodology".

# Helper: PowerPoint's Shape.Height/.Width/.Left/.Top are expressed in points
# (float, internally stored as 32-bit float) while the OOXML <a:ext>/<a:off>
# values are EMU (1 pt = 12700 EMU). A naive emu/12700 division can land
# just under the target EMU once it is rounded through a 32-bit float and
# truncated back to EMU, so nudge upward (in tiny steps) until the value
# round-trips to exactly the desired EMU amount.
function ConvertTo-PointsForEmu($targetEmu) {
    $pts = $targetEmu / 12700.0
    for ($i = 0; $i -lt 200000; $i++) {
        $f = [double]([float]$pts)
        $emu = [math]::Floor($f * 12700.0)
        if ($emu -eq $targetEmu) {
            return $pts
        }
        $pts = $pts + 0.0000001
    }
    return $targetEmu / 12700.0
}

$p = $ppt.ActivePresentation

# --- Slide 15: TextBox 7 -----------------------------------------------
$s15 = $p.Slides.Item(15)
$tb = $s15.Shapes.Item("TextBox 7")
$tbText = $tb.TextFrame.TextRange

$tbText.Paragraphs(1).Runs(1).Font.Size = 20
$tbText.Paragraphs(2).Runs(1).Font.Size = 20

$tb.Height = ConvertTo-PointsForEmu 1631216

# --- Slide 4: Content Placeholder 2 -------------------------------------
$s4 = $p.Slides.Item(4)
$cp = $s4.Shapes.Item("Content Placeholder 2")
$cpText = $cp.TextFrame.TextRange

$quote1 = [char]0x201C
$quote2 = [char]0x201D
$cpText.Paragraphs(1).Runs(1).Text = "Project is inspired by the " + $quote1 + "Analyze Customer Lifetime Value" + $quote2 + " Databricks Solution Accelerator"
$cpText.Paragraphs(2).Runs(1).Text = "Analyze Customer Lifetime Value"

# --- Slide 6: Methodology table ------------------------------------------
$s6 = $p.Slides.Item(6)
$tblShape = $s6.Shapes.Item("Content Placeholder 5")
$cell = $tblShape.Table.Cell(1, 2)
$cell.Shape.TextFrame.TextRange.Text = "Project Methodology"
